$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New header labels for columns U, V, W
$ws.Range("U1").Value = "belly toward high z value"
$ws.Range("V1").Value = "belly toward high x value"
$ws.Range("W1").Value = "belly toward high y value"

# New data (x1 ind / y1 ind columns J,K,M,N and z ind columns P,Q, plus belly U,V,W)
$data = @{
    2  = @{ J=345; K=643; M=313;  N=923;  P=6;  Q=92; U=1;  V=0;  W=0 }
    3  = @{ J=587; K=855; M=237;  N=897;  P=6;  Q=92; U=1;  V=0;  W=0 }
    4  = @{ J=338; K=642; M=217;  N=898;  P=6;  Q=96; U=-1; V=0;  W=0 }
    5  = @{ J=335; K=662; M=345;  N=991;  P=9;  Q=97; U=1;  V=0;  W=0 }
    6  = @{ J=384; K=709; M=327;  N=983;  P=3;  Q=97; U=-1; V=0;  W=0 }
    7  = @{ J=270; K=624; M=175;  N=817;  P=5;  Q=91; U=-1; V=0;  W=0 }
    8  = @{ J=442; K=878; M=1090; N=1914; P=4;  Q=47; U=-1; V=0;  W=0 }
    9  = @{ J=503; K=970; M=1152; N=1936; P=10; Q=55; U=0;  V=-1; W=0 }
    10 = @{ J=387; K=836; M=1108; N=1905; P=12; Q=55; U=0;  V=1;  W=0 }
    11 = @{ J=457; K=857; M=1134; N=1881; P=7;  Q=60; U=0;  V=1;  W=0 }
    12 = @{ J=404; K=832; M=1141; N=1929; P=3;  Q=49; U=0;  V=-1; W=0 }
    13 = @{ J=524; K=934; M=976;  N=1878; P=3;  Q=50; U=0;  V=1;  W=0 }
    14 = @{ J=372; K=774; M=1127; N=1918; P=8;  Q=65; U=0;  V=-1; W=0 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$row").Value = $vals[$col]
    }
}

# Column widths for U, V, W (recalculated bestFit widths after header text change)
$ws.Range("U1").EntireColumn.ColumnWidth = 22.6
$ws.Range("V1:W1").EntireColumn.ColumnWidth = 22.6

# Update the selected cell to match the post-edit state
$ws.Range("X15").Select()
